$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: add "Data1" / "Data" headers in C6 / D6 (plain default data style, same as B6)
$ws.Range("B6").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("C6").Value = "Data1"
$ws.Range("D6").Value = "Data"

# Row 7: C7/D7 switch from text "Data" to numeric sample values
$ws.Range("C7").Value = 123
$ws.Range("D7").Value = 456

# Row 8: add a new "Admin " / "admin123" pair, reusing row 7's formatting for A8
# (blue font + border), then turning word-wrap on for that cell.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Admin "
$ws.Range("A8").WrapText = $true

$ws.Range("B3").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = "admin123"

# Row 9: new row with "***"
$ws.Range("A4").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "***"

$ws.Range("C12").Select()
